# Actualización desde MV -datos-
# Update the last few months of data (rows 218-222) with revised figures
# and append a new row (223) for "01-06-2021".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revised values for existing rows 218-222 (columns B, C, D, F change;
#     A, E, G stay as they were) ---
$revised = @{
    218 = @{ B = -179;  C = 321;   D = 2684;  F = -2815 }
    219 = @{ B = -149;  C = 11;    D = -1031; F = 485 }
    220 = @{ B = -590;  C = -839;  D = -209;  F = -308 }
    221 = @{ B = -995;  C = -3117; D = -3427; F = 3499 }
    222 = @{ B = -1507; C = 484;   D = -7644; F = 681 }
}

foreach ($row in $revised.Keys) {
    $vals = $revised[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("F$row").Value = $vals.F
}

# --- New row 223: "01-06-2021" ---
$ws.Range("A223").Value = "01-06-2021"
$ws.Range("B223").Value = -2293
$ws.Range("C223").Value = -803
$ws.Range("D223").Value = -642
$ws.Range("E223").Value = 92
$ws.Range("F223").Value = 1655
$ws.Range("G223").Value = -2596
